$d = $word.ActiveDocument

# --- Change 1: "{{ rg_cliente }}" -> "{{rg_cliente}}" (drop the leading
# and trailing space runs around the placeholder, keep its highlight) ---

# Remove the space right after the opening "{{" (before "rg_cliente").
# "RG n...: {{" only occurs once in the document (right before the
# rg_cliente placeholder), so this anchor is unambiguous.
$rng1b = $d.Content
$found1b = $rng1b.Find.Execute("RG n" + [char]0x00BA + ": {{", $true, $false, `
                                $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1b) {
    throw "Could not locate 'RG n<degree>: {{' anchor before rg_cliente"
}
$space1 = $d.Range($rng1b.End, $rng1b.End + 1)
if ($space1.Text -eq " ") {
    $space1.Text = ""
}

# Remove the space right before the closing "}}" (after "rg_cliente").
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("rg_cliente", $true, $false, $false, $false, `
                              $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate 'rg_cliente' placeholder text"
}
$space2 = $d.Range($rng2.End, $rng2.End + 1)
if ($space2.Text -eq " ") {
    $space2.Text = ""
}

# --- Change 2: rename placeholder "nombre_inmueble" -> "cod_referencia"
# and drop its yellow highlight ---

$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$find3.Text = "nombre_inmueble"
$find3.Replacement.Text = "cod_referencia"
$find3.Replacement.Highlight = $false
$found3 = $find3.Execute($find3.Text, $true, $false, $false, $false, $false, `
                          $true, 1, $false, $find3.Replacement.Text, 2, $true)
if (-not $found3) {
    throw "Could not locate 'nombre_inmueble' placeholder text"
}
